$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.169.86'
$ws.Range('E2').Value = '  -0.08%  '

$ws.Range('D3').Value = '1.817.65'
$ws.Range('E3').Value = '  -0.63%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.79%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.81'
$ws.Range('E5').Value = '  -1.21%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5907'
$ws.Range('E6').Value = '  -2.53%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.51%  '

$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2723'
$ws.Range('E8').Value = '  -3.83%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06801'
$ws.Range('E9').Value = '  -4.51%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.99'
$ws.Range('E10').Value = '  -4.53%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07546'
$ws.Range('E11').Value = '  -1.49%  '

$ws.Range('D12').Value = '1.827.92'
$ws.Range('E12').Value = '  +0.13%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.648'
$ws.Range('E13').Value = '  -2.97%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6154'
$ws.Range('E14').Value = '  -4.24%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009553'
$ws.Range('E15').Value = '  -5.09%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '75.54'
$ws.Range('E16').Value = '  -5.12%  '

$ws.Range('D17').Value = '28.921.39'
$ws.Range('E17').Value = '  -0.92%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.440'
$ws.Range('E18').Value = '  -9.64%  '

$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.50%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '207.80'
$ws.Range('E20').Value = '  -10.33%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.39'
$ws.Range('E21').Value = '  -3.47%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.721'
$ws.Range('E22').Value = '  -4.65%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.010'
$ws.Range('E23').Value = '  +1.37%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '155.08'
$ws.Range('E24').Value = '  -0.23%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.819'
$ws.Range('E25').Value = '  -2.91%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1271'
$ws.Range('E26').Value = '  -1.35%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.25'
$ws.Range('E27').Value = '  -2.83%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06394'
$ws.Range('E28').Value = '  -7.90%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.416'
$ws.Range('E29').Value = '  -3.12%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.436'
$ws.Range('E30').Value = '  -1.05%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.709'
$ws.Range('E31').Value = '  -2.65%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.702'
$ws.Range('E32').Value = '  -3.50%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.073'
$ws.Range('E33').Value = '  -5.87%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.691'
$ws.Range('E34').Value = '  -2.06%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.544'
$ws.Range('E35').Value = '  +0.48%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6300'
$ws.Range('E36').Value = '  -5.12%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.766'
$ws.Range('E37').Value = '  +0.41%  '

$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.517'
$ws.Range('E38').Value = '  -1.35%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01721'
$ws.Range('E39').Value = '  -2.79%  '

$ws.Range('D40').Value = '1.125.51'
$ws.Range('E40').Value = '  -8.98%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8803'
$ws.Range('E41').Value = '  -5.50%  '

$ws.Range('E42').Value = '  +0.35%  '

$ws.Range('D43').Value = '1.983.61'
$ws.Range('E43').Value = '  -0.80%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.60'
$ws.Range('E44').Value = '  -0.63%  '

$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000116'
$ws.Range('E45').Value = '  -0.98%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.04'
$ws.Range('E46').Value = '  -3.96%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.584'
$ws.Range('E47').Value = '  -3.70%  '

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05508'
$ws.Range('E48').Value = '  -1.53%  '

$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4543'
$ws.Range('E49').Value = '  -0.27%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.247'
$ws.Range('E50').Value = '  -3.67%  '

$ws.Range('B51').Value = 'Frax'
$ws.Range('C51').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9922'
$ws.Range('E51').Value = '  -0.72%  '
